$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty column I (no longer used in the updated table)
$ws.Columns("I:I").Delete()

# Insert a new row for "range_nr" before the current row 32 (Federalland_pcnt2),
# shifting Federalland_pcnt2 / Waterland_pcnt2 / Ruralland_pcnt2 down one row.
$ws.Rows("32:32").Insert()

# Update recomputed values for row 27 (pasture_nr)
$ws.Range("E27").Value = 22.451137542724609
$ws.Range("F27").Value = 21.094793319702148
$ws.Range("G27").Value = 22.477123260498047
$ws.Range("H27").Value = 20.875175476074219

# Update recomputed values for row 28 (CRP_nr)
$ws.Range("G28").Value = 62.593730926513672
$ws.Range("H28").Value = 64.778839111328125

# Row 29 (crop_nr) values are unchanged

# Update recomputed values for row 30 (forest_nr)
$ws.Range("F30").Value = 20.544464111328125
$ws.Range("G30").Value = 19.215654373168945
$ws.Range("H30").Value = 16.132402420043945

# Update recomputed values for row 31 (urban_nr)
$ws.Range("E31").Value = 22257.21875
$ws.Range("F31").Value = 28130.1640625
$ws.Range("G31").Value = 41444.46484375
$ws.Range("H31").Value = 43307.26953125

# Fill the newly-inserted row 32 with the "range_nr" label and values
# (matching the updated row 27 "pasture_nr" values above)
$ws.Range("A32").Value = "range_nr"
$ws.Range("E32").Value = 22.451137542724609
$ws.Range("F32").Value = 21.094793319702148
$ws.Range("G32").Value = 22.477123260498047
$ws.Range("H32").Value = 20.875175476074219
